# "Final Formatting of Presentation"
#
# Slide 7 ("DDS Analytics - Data Insight"), shape 2 ("TextBox 73") holds the
# life-satisfaction formula sentence. It was typed using the Unicode
# "mathematical sans/serif bold italic" alphanumeric block (fake-bold
# Unicode glyphs) instead of plain text with real bold formatting, and had a
# typo / missing final punctuation. Replace it with normal ASCII text,
# fixing the wording/typos and ending the sentence with a period. The
# corrected word "Involvment" (flagged by the spell checker in the source
# file) and the trailing "." are kept as their own runs, matching how
# PowerPoint split them when the misspelling was flagged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Assigning .Text rewrites the paragraph as a single run while keeping the
# existing run-level formatting (lang="en-US" dirty="0", latin "+mj-lt").
$tr.Text = "Life Satisfaction is average of Environmental Satisfaction, Job Satisfaction, Relationship Satisfaction, Work Life Balance and Job Involvment."

# Peel "Involvment" and the trailing "." off into their own runs. Touching
# (re-assigning) a character-level Font property forces PowerPoint to split
# the run at that point; Font.Spacing is used because re-applying its own
# current value is a pure no-op that does not add any extra formatting to
# the run (so all three runs keep identical rPr / latin typeface).
$involvment = $tr.Characters(132, 10)
$involvment.Font.Spacing = $involvment.Font.Spacing

$period = $tr.Characters(142, 1)
$period.Font.Spacing = $period.Font.Spacing
